$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-12-20"
$ws.Range("B1").Value = "December 2022 (through December 20)"

$ws.Range("Z2").Value = 4
$ws.Range("AX2").Value = 4
$ws.Range("BJ2").Value = 4
$ws.Range("BV2").Value = 5
$ws.Range("CH2").Value = 3
$ws.Range("N3").Value = 2
$ws.Range("AL3").Value = 2
$ws.Range("BV3").Value = 6
$ws.Range("Z7").Value = 6
$ws.Range("N14").Value = 6
$ws.Range("AL14").Value = 5
$ws.Range("Z15").Value = 2
$ws.Range("N16").Value = 3
$ws.Range("AX20").Value = 1
$ws.Range("N25").Value = 2
$ws.Range("B26").Value = 2
$ws.Range("CH29").Value = 2
$ws.Range("CH30").Value = 1
$ws.Range("B40").Value = 2
$ws.Range("BJ43").Value = 2
$ws.Range("BV45").Value = 1
$ws.Range("BJ59").Value = 3
$ws.Range("BV66").Value = 1
$ws.Range("B91").Value = 1
$ws.Range("N92").Value = 3
